$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "5682953 (not retierd  2560405, 2572012)"
$ws.Range("A1").Characters(9, 100).Font.Color = RGB(192, 0, 0)
